$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$arr = New-Object "object[,]" 1,36
$arr[0,0] = 2.42
$arr[0,1] = 2.5
$arr[0,2] = 3.55
$arr[0,3] = 3.65
$arr[0,4] = 3.15
$arr[0,5] = 3.3
$arr[0,6] = 1.51
$arr[0,7] = 1.1
$arr[0,8] = 3.05
$arr[0,9] = 1.43
$arr[0,10] = 1.67
$arr[0,11] = 2.36
$arr[0,12] = 1.25
$arr[0,13] = 4.4
$arr[0,14] = 1.94
$arr[0,15] = 1.94
$arr[0,16] = 1.32
$arr[0,17] = 1.51
$arr[0,18] = 12.5
$arr[0,19] = 11
$arr[0,20] = 22
$arr[0,21] = 70
$arr[0,22] = 9
$arr[0,23] = 7.2
$arr[0,24] = 15
$arr[0,25] = 48
$arr[0,26] = 14.5
$arr[0,27] = 12
$arr[0,28] = 21
$arr[0,29] = 70
$arr[0,30] = 34
$arr[0,31] = 32
$arr[0,32] = 55
$arr[0,33] = 140
$arr[0,34] = 30
$arr[0,35] = 60
$ws.Range("F2:AO2").Value = $arr

# Row 3
$arr = New-Object "object[,]" 1,36
$arr[0,0] = 1.72
$arr[0,1] = 1.8
$arr[0,2] = 5.7
$arr[0,3] = 8.199999999999999
$arr[0,4] = 3.55
$arr[0,5] = 3.95
$arr[0,6] = 1.54
$arr[0,7] = 1.09
$arr[0,8] = 2.9
$arr[0,9] = 1.45
$arr[0,10] = 1.63
$arr[0,11] = 2.46
$arr[0,12] = 1.23
$arr[0,13] = 4.4
$arr[0,14] = 2.14
$arr[0,15] = 1.69
$arr[0,16] = 1.14
$arr[0,17] = 2.04
$arr[0,18] = 13
$arr[0,19] = 19.5
$arr[0,20] = 60
$arr[0,21] = 1000
$arr[0,22] = 7.8
$arr[0,23] = 8.800000000000001
$arr[0,24] = 980
$arr[0,25] = 1000
$arr[0,26] = 9.6
$arr[0,27] = 11
$arr[0,28] = 1000
$arr[0,29] = 1000
$arr[0,30] = 22
$arr[0,31] = 27
$arr[0,32] = 1000
$arr[0,33] = 1000
$arr[0,34] = 19
$arr[0,35] = 1000
$ws.Range("F3:AO3").Value = $arr

# Row 4
$arr = New-Object "object[,]" 1,36
$arr[0,0] = 2.92
$arr[0,1] = 3.15
$arr[0,2] = 2.72
$arr[0,3] = 2.9
$arr[0,4] = 3.05
$arr[0,5] = 3.25
$arr[0,6] = 1.5
$arr[0,7] = 1.08
$arr[0,8] = 2.96
$arr[0,9] = 1.44
$arr[0,10] = 1.66
$arr[0,11] = 2.28
$arr[0,12] = 1.24
$arr[0,13] = 4.6
$arr[0,14] = 1.87
$arr[0,15] = 1.93
$arr[0,16] = 1.52
$arr[0,17] = 1.46
$arr[0,18] = 12.5
$arr[0,19] = 11.5
$arr[0,20] = 18.5
$arr[0,21] = 55
$arr[0,22] = 12
$arr[0,23] = 8.6
$arr[0,24] = 15.5
$arr[0,25] = 44
$arr[0,26] = 24
$arr[0,27] = 14.5
$arr[0,28] = 21
$arr[0,29] = 60
$arr[0,30] = 65
$arr[0,31] = 50
$arr[0,32] = 70
$arr[0,33] = 160
$arr[0,34] = 55
$arr[0,35] = 1000
$ws.Range("F4:AO4").Value = $arr

# Row 5
$arr = New-Object "object[,]" 1,36
$arr[0,0] = 1.64
$arr[0,1] = 1.72
$arr[0,2] = 6.2
$arr[0,3] = 7.2
$arr[0,4] = 3.75
$arr[0,5] = 4.2
$arr[0,6] = 1.47
$arr[0,7] = 1.09
$arr[0,8] = 3.1
$arr[0,9] = 1.41
$arr[0,10] = 1.73
$arr[0,11] = 2.2
$arr[0,12] = 1.26
$arr[0,13] = 4.2
$arr[0,14] = 2.12
$arr[0,15] = 1.78
$arr[0,16] = 1.16
$arr[0,17] = 2.38
$arr[0,18] = 14.5
$arr[0,19] = 22
$arr[0,20] = 65
$arr[0,21] = 270
$arr[0,22] = 7.8
$arr[0,23] = 11
$arr[0,24] = 32
$arr[0,25] = 150
$arr[0,26] = 11
$arr[0,27] = 12.5
$arr[0,28] = 32
$arr[0,29] = 150
$arr[0,30] = 21
$arr[0,31] = 25
$arr[0,32] = 60
$arr[0,33] = 230
$arr[0,34] = 16.5
$arr[0,35] = 1000
$ws.Range("F5:AO5").Value = $arr

# Row 6
$arr = New-Object "object[,]" 1,36
$arr[0,0] = 1.86
$arr[0,1] = 1.89
$arr[0,2] = 5.4
$arr[0,3] = 5.8
$arr[0,4] = 3.45
$arr[0,5] = 3.65
$arr[0,6] = 1.48
$arr[0,7] = 1.08
$arr[0,8] = 3.1
$arr[0,9] = 1.39
$arr[0,10] = 1.75
$arr[0,11] = 2.1
$arr[0,12] = 1.27
$arr[0,13] = 3.9
$arr[0,14] = 1.96
$arr[0,15] = 1.84
$arr[0,16] = 1.21
$arr[0,17] = 2.12
$arr[0,18] = 12.5
$arr[0,19] = 16.5
$arr[0,20] = 40
$arr[0,21] = 150
$arr[0,22] = 7.6
$arr[0,23] = 8.4
$arr[0,24] = 22
$arr[0,25] = 85
$arr[0,26] = 10
$arr[0,27] = 10.5
$arr[0,28] = 1000
$arr[0,29] = 95
$arr[0,30] = 20
$arr[0,31] = 22
$arr[0,32] = 1000
$arr[0,33] = 160
$arr[0,34] = 15.5
$arr[0,35] = 120
$ws.Range("F6:AO6").Value = $arr

# Row 7
$arr = New-Object "object[,]" 1,36
$arr[0,0] = 1.07
$arr[0,1] = 1.09
$arr[0,2] = 28
$arr[0,3] = 1000
$arr[0,4] = 14
$arr[0,5] = 21
$arr[0,6] = 1.22
$arr[0,7] = 1.01
$arr[0,8] = 7.8
$arr[0,9] = 1.11
$arr[0,10] = 3.3
$arr[0,11] = 1.34
$arr[0,12] = 1.93
$arr[0,13] = 1.74
$arr[0,14] = 3.15
$arr[0,15] = 1.37
$arr[0,16] = 1.01
$arr[0,17] = 12
$arr[0,18] = 1000
$arr[0,19] = 180
$arr[0,20] = 1000
$arr[0,21] = 1000
$arr[0,22] = 1000
$arr[0,23] = 1000
$arr[0,24] = 1000
$arr[0,25] = 1000
$arr[0,26] = 1000
$arr[0,27] = 1000
$arr[0,28] = 1000
$arr[0,29] = 1000
$arr[0,30] = 1000
$arr[0,31] = 1000
$arr[0,32] = 1000
$arr[0,33] = 1000
$arr[0,34] = 2.5
$arr[0,35] = 1000
$ws.Range("F7:AO7").Value = $arr

# Row 8
$arr = New-Object "object[,]" 1,36
$arr[0,0] = 1.64
$arr[0,1] = 1.73
$arr[0,2] = 6.6
$arr[0,3] = 7.6
$arr[0,4] = 3.55
$arr[0,5] = 3.95
$arr[0,6] = 1.5
$arr[0,7] = 1.1
$arr[0,8] = 2.84
$arr[0,9] = 1.46
$arr[0,10] = 1.62
$arr[0,11] = 2.34
$arr[0,12] = 1.23
$arr[0,13] = 4.5
$arr[0,14] = 2.26
$arr[0,15] = 1.67
$arr[0,16] = 1.15
$arr[0,17] = 2.36
$arr[0,18] = 12.5
$arr[0,19] = 23
$arr[0,20] = 1000
$arr[0,21] = 330
$arr[0,22] = 7.6
$arr[0,23] = 10.5
$arr[0,24] = 36
$arr[0,25] = 180
$arr[0,26] = 10.5
$arr[0,27] = 13
$arr[0,28] = 34
$arr[0,29] = 170
$arr[0,30] = 19.5
$arr[0,31] = 26
$arr[0,32] = 1000
$arr[0,33] = 250
$arr[0,34] = 16
$arr[0,35] = 1000
$ws.Range("F8:AO8").Value = $arr

# Row 9
$arr = New-Object "object[,]" 1,36
$arr[0,0] = 1.47
$arr[0,1] = 1.54
$arr[0,2] = 8.6
$arr[0,3] = 12
$arr[0,4] = 4.1
$arr[0,5] = 4.7
$arr[0,6] = 1.01
$arr[0,7] = 1.07
$arr[0,8] = 3.35
$arr[0,9] = 1.35
$arr[0,10] = 1.83
$arr[0,11] = 2.04
$arr[0,12] = 1.3
$arr[0,13] = 3.75
$arr[0,14] = 2.2
$arr[0,15] = 1.7
$arr[0,16] = 1.1
$arr[0,17] = 2.84
$arr[0,18] = 16
$arr[0,19] = 29
$arr[0,20] = 1000
$arr[0,21] = 460
$arr[0,22] = 8
$arr[0,23] = 12
$arr[0,24] = 42
$arr[0,25] = 220
$arr[0,26] = 9.199999999999999
$arr[0,27] = 12.5
$arr[0,28] = 36
$arr[0,29] = 190
$arr[0,30] = 14.5
$arr[0,31] = 21
$arr[0,32] = 60
$arr[0,33] = 260
$arr[0,34] = 10
$arr[0,35] = 1000
$ws.Range("F9:AO9").Value = $arr

# Row 10
$arr = New-Object "object[,]" 1,36
$arr[0,0] = 1.25
$arr[0,1] = 1.29
$arr[0,2] = 1.49
$arr[0,3] = 22
$arr[0,4] = 5
$arr[0,5] = 7
$arr[0,6] = 1.36
$arr[0,7] = 1.05
$arr[0,8] = 3.45
$arr[0,9] = 1.26
$arr[0,10] = 1.88
$arr[0,11] = 1.67
$arr[0,12] = 1.39
$arr[0,13] = 2.72
$arr[0,14] = 2.56
$arr[0,15] = 1.52
$arr[0,16] = 1.05
$arr[0,17] = 4.3
$arr[0,18] = 19
$arr[0,19] = 1000
$arr[0,20] = 250
$arr[0,21] = 1000
$arr[0,22] = 8.6
$arr[0,23] = 18
$arr[0,24] = 1000
$arr[0,25] = 600
$arr[0,26] = 7
$arr[0,27] = 12.5
$arr[0,28] = 1000
$arr[0,29] = 420
$arr[0,30] = 9
$arr[0,31] = 17.5
$arr[0,32] = 1000
$arr[0,33] = 450
$arr[0,34] = 5.6
$arr[0,35] = 1000
$ws.Range("F10:AO10").Value = $arr

# Row 11
$arr = New-Object "object[,]" 1,36
$arr[0,0] = 2.36
$arr[0,1] = 2.5
$arr[0,2] = 3.3
$arr[0,3] = 3.75
$arr[0,4] = 3.25
$arr[0,5] = 3.6
$arr[0,6] = 1.5
$arr[0,7] = 1.09
$arr[0,8] = 2.78
$arr[0,9] = 1.44
$arr[0,10] = 1.64
$arr[0,11] = 2.26
$arr[0,12] = 1.23
$arr[0,13] = 4.6
$arr[0,14] = 1.93
$arr[0,15] = 1.86
$arr[0,16] = 1.37
$arr[0,17] = 1.67
$arr[0,18] = 12.5
$arr[0,19] = 1000
$arr[0,20] = 24
$arr[0,21] = 1000
$arr[0,22] = 10.5
$arr[0,23] = 8.800000000000001
$arr[0,24] = 15.5
$arr[0,25] = 1000
$arr[0,26] = 1000
$arr[0,27] = 1000
$arr[0,28] = 1000
$arr[0,29] = 1000
$arr[0,30] = 1000
$arr[0,31] = 1000
$arr[0,32] = 1000
$arr[0,33] = 180
$arr[0,34] = 1000
$arr[0,35] = 1000
$ws.Range("F11:AO11").Value = $arr
